$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 81, shifting rows 81:89 down to 84:92.
$ws.Rows("81:83").Insert()

# Common (unchanged) column values shared by every row in this block.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103003
$categoria   = "Damasco"

$newRows = @(
    @{ Row = 81; Fecha = 44918; Variedad = "Dina"; Calidad = "Especial"; Volumen = 240; PrecioMin = 23000; PrecioMax = 24000; PrecioProm = 23500; Unidad = "$/caja 16 kilos"; Origen = "Región Metropolitana"; PrecioKg = 1469; KgUnidad = 16 },
    @{ Row = 82; Fecha = 44918; Variedad = "Dina"; Calidad = "Primera"; Volumen = 160; PrecioMin = 20000; PrecioMax = 21000; PrecioProm = 20500; Unidad = "$/caja 16 kilos"; Origen = "Región Metropolitana"; PrecioKg = 1281; KgUnidad = 16 },
    @{ Row = 83; Fecha = 44918; Variedad = "Dina"; Calidad = "Segunda"; Volumen = 160; PrecioMin = 15000; PrecioMax = 16000; PrecioProm = 15500; Unidad = "$/caja 16 kilos"; Origen = "Región Metropolitana"; PrecioKg = 969;  KgUnidad = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PrecioMin
    $ws.Cells.Item($row, 15).Value = $r.PrecioMax
    $ws.Cells.Item($row, 16).Value = $r.PrecioProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
